{"js": "// The document repeats a \"\u25b2\u8868 8-2-N ...\" caption pattern for every table\n// (\u8868 8-2-2 through \u8868 8-2-10). The edit removes the leading \"\u25b2\" run from\n// only the very first caption paragraph (\"\u25b2\u8868 8-2-2 \u56de\u6536\u9805\u76ee\"), leaving the\n// caption text run and all of the other \"\u25b2\" captions untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"\u25b2\u8868 8-2-2 \u56de\u6536\u9805\u76ee\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  // Scope the search to this single paragraph so only its own \"\u25b2\" run is\n  // matched (the same glyph also starts every other table caption).\n  const paragraphRange = target.getRange();\n  const hits = paragraphRange.search(\"\u25b2\", { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n\n  if (hits.items.length > 0) {\n    hits.items[0].delete();\n    await context.sync();\n  }\n}\n", "ps1": "# The document repeats a \"\u25b2\u8868 8-2-N ...\" caption pattern for every table\n# (\u8868 8-2-2 through \u8868 8-2-10). The edit removes only the leading \"\u25b2\"\n# character from the very first caption paragraph (\"\u25b2\u8868 8-2-2 \u56de\u6536\u9805\u76ee\"),\n# leaving the caption text and every other \"\u25b2\" caption untouched.\n$d = $word.ActiveDocument\n\n$triangle = [char]0x25B2\n$targetText = $triangle + \"\u8868 8-2-2 \u56de\u6536\u9805\u76ee\" + [char]13\n\n$target = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text -eq $targetText) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    # Scope Find to just this paragraph so only its own leading \"\u25b2\" is\n    # matched (the same glyph also starts every other table caption).\n    $rng = $target.Range\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $triangle\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = \"\"\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
